$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price / volume(1h) snapshot.
# Price cells (column D) whose new text looks like a plain decimal number
# are forced to Text format first so Excel keeps the literal string
# (e.g. "1.00", "0.580") instead of silently coercing it to a number.

$ws.Range("D2").Value = '66.610.53'
$ws.Range("E2").Value = '  +3.91%  '

$ws.Range("D3").Value = '3.501.26'
$ws.Range("E3").Value = '  +2.45%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.90'
$ws.Range("E5").Value = '  +3.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '170.26'
$ws.Range("E6").Value = '  +5.93%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.02%  '

$ws.Range("D8").Value = '3.497.85'
$ws.Range("E8").Value = '  +2.24%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.580'
$ws.Range("E9").Value = '  +4.60%  '

$ws.Range("E10").Value = '  +0.92%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.125'
$ws.Range("E11").Value = '  +4.29%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.438'
$ws.Range("E12").Value = '  +3.17%  '

$ws.Range("D13").Value = '4.105.94'
$ws.Range("E13").Value = '  +2.45%  '

$ws.Range("E14").Value = '  -0.51%  '

$ws.Range("E15").Value = '  +4.06%  '

$ws.Range("E16").Value = '  +1.80%  '

$ws.Range("D17").Value = '66.644.90'
$ws.Range("E17").Value = '  +3.82%  '

$ws.Range("D18").Value = '3.496.92'
$ws.Range("E18").Value = '  +1.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.33'
$ws.Range("E19").Value = '  +3.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.03'
$ws.Range("E20").Value = '  +3.20%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '388.57'
$ws.Range("E21").Value = '  +2.68%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '8.01'
$ws.Range("E22").Value = '  +1.86%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '73.06'
$ws.Range("E23").Value = '  +2.29%  '

$ws.Range("E24").Value = '  +0.08%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.534'
$ws.Range("E25").Value = '  +3.02%  '

$ws.Range("E26").Value = '  +4.46%  '

$ws.Range("E28").Value = '  +1.85%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  -0.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.34'
$ws.Range("E30").Value = '  +4.45%  '

$ws.Range("E31").Value = '  +5.58%  '

$ws.Range("E32").Value = '  +2.68%  '

$ws.Range("E33").Value = '  +2.62%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.45'

$ws.Range("E35").Value = '  +0.05%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.61'
$ws.Range("E36").Value = '  +6.74%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '162.68'
$ws.Range("E37").Value = '  +2.26%  '

$ws.Range("E38").Value = '  +2.40%  '

$ws.Range("E39").Value = '  +3.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.70'
$ws.Range("E40").Value = '  +5.73%  '

$ws.Range("E41").Value = '  +2.20%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '26.39'
$ws.Range("E42").Value = '  +2.50%  '

$ws.Range("D43").Value = '2.825.93'
$ws.Range("E43").Value = '  +0.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.64'
$ws.Range("E44").Value = '  +1.32%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '26.64'
$ws.Range("E45").Value = '  +1.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.02'
$ws.Range("E46").Value = '  -0.28%  '

$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0312'
$ws.Range("E47").Value = '  +2.48%  '

$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '354.98'
$ws.Range("E48").Value = '  +3.65%  '

$ws.Range("B49").Value = 'dogwifhat'
$ws.Range("C49").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.51'
$ws.Range("E49").Value = '  +5.55%  '

$ws.Range("E50").Value = '  +2.61%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '33.72'
$ws.Range("E51").Value = '  +11.74%  '
